# Add homework row for 2020-03-12 (row 79)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 79

$ws.Cells.Item($row, 1).Value = 1583971200

# Force B/C to remain plain text (avoid Excel auto-converting the date-like
# string to a date serial, and the numeric-looking id to a number), then
# clear the temporary number format so the cell style stays the default.
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2020-03-12"
$ws.Cells.Item($row, 2).ClearFormats()

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "0216"
$ws.Cells.Item($row, 3).ClearFormats()

$ws.Cells.Item($row, 4).Value = "SPRING"
$ws.Cells.Item($row, 5).Value = 0.19
$ws.Cells.Item($row, 6).Value = 0.195
$ws.Cells.Item($row, 7).Value = 0.18
$ws.Cells.Item($row, 8).Value = 0.18
$ws.Cells.Item($row, 9).Value = 643000
